# Update syllabus for OLS
# Rename the three "OLS ..." topic rows to their "Ordinary least squares ..." equivalents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D19").Value = "Ordinary least squares - foundation and diagnostics"
$ws.Range("D20").Value = "Ordinary least squares - interpretation/hypothesis testing"
$ws.Range("D21").Value = "Ordinary least squares - multivariable/interaction terms"

# Move the active selection to D22, matching the author's final cursor position.
$ws.Range("D22").Select()
